$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.839.68'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '1.562.22'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.78'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.77'
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = '1.783.20'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = '1.577.97'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '26.853.13'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.21'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.96'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.35'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.20'
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.69'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.72'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.90'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  -1.39%  '
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('E31').Value = '  -3.87%  '
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '1.404.02'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.914'
$ws.Range('E37').Value = '  -3.64%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.527'
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.810'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.996'
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.51'
$ws.Range('E43').Value = '  +5.86%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.33'
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('D47').Value = '1.697.11'
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.57'
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('D50').Value = '0.0₇0971'
$ws.Range('E50').Value = '  -2.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0949'
$ws.Range('E51').Value = '  +0.92%  '
